## "add lan 3 agin" -- turn the trailing blank paragraph into a third
## numbered "Lần 3" list item, matching the existing "Lần 1" / "Lần 2"
## entries (same ListParagraph style + numId 4 numbering).

$d = $word.ActiveDocument

# The target is the very last paragraph in the document body -- currently
# empty, formatted with `ind left=720` (left-over list indent, no real
# numbering applied yet).
$target = $d.Paragraphs.Last

# Replace it with a paragraph that has the same pPr shape (style + numPr
# + rPr) as the two preceding "Lần 1" / "Lần 2" list items, reusing
# numId 4 so the numbering continues as 1, 2, 3, plus a run carrying the
# new "Lần 3" text with matching run formatting.
$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr>' +
            '<w:ilvl w:val="0"/>' +
            '<w:numId w:val="4"/>' +
        '</w:numPr>' +
        '<w:rPr>' +
            '<w:rFonts w:cstheme="minorHAnsi"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:sz w:val="40"/>' +
            '<w:szCs w:val="40"/>' +
            '<w:lang w:val="vi-VN"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:cstheme="minorHAnsi"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:sz w:val="40"/>' +
            '<w:szCs w:val="40"/>' +
            '<w:lang w:val="vi-VN"/>' +
        '</w:rPr>' +
        '<w:lastRenderedPageBreak/>' +
        '<w:t>Lần 3</w:t>' +
    '</w:r>' +
'</w:p>'

$target.Range.InsertXML($paraXml)
